$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.776.00'
$ws.Range('E2').Value = '  -7.09%  '
$ws.Range('D3').Value = '2.545.70'
$ws.Range('E3').Value = '  -2.02%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '296.68'
$ws.Range('E5').Value = '  -4.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.54'
$ws.Range('E6').Value = '  -7.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.572'
$ws.Range('E7').Value = '  -4.44%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -5.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.65'
$ws.Range('E10').Value = '  -8.48%  '
$ws.Range('E11').Value = '  -4.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.60'
$ws.Range('E12').Value = '  -6.42%  '
$ws.Range('D13').Value = '2.933.31'
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = '2.532.58'
$ws.Range('E15').Value = '  -2.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.862'
$ws.Range('E16').Value = '  -5.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.10'
$ws.Range('E17').Value = '  -5.09%  '
$ws.Range('D18').Value = '42.795.86'
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.63'
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0975'
$ws.Range('E20').Value = '  -4.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.46'
$ws.Range('E21').Value = '  -2.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.46'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '260.01'
$ws.Range('E24').Value = '  -5.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '29.49'
$ws.Range('E25').Value = '  -0.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.11'
$ws.Range('E26').Value = '  -6.50%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.98'
$ws.Range('E28').Value = '  -7.29%  '
$ws.Range('E29').Value = '  -4.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.02'
$ws.Range('E30').Value = '  -5.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.89'
$ws.Range('E31').Value = '  -5.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '150.58'
$ws.Range('E32').Value = '  -3.05%  '
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.38'
$ws.Range('E34').Value = '  -5.37%  '
$ws.Range('E35').Value = '  -2.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0791'
$ws.Range('E36').Value = '  -5.57%  '
$ws.Range('E37').Value = '  -6.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.12'
$ws.Range('E38').Value = '  +13.62%  '
$ws.Range('E39').Value = '  -3.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.09'
$ws.Range('E40').Value = '  +1.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.42'
$ws.Range('E41').Value = '  -4.69%  '
$ws.Range('E42').Value = '  -6.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.81'
$ws.Range('E43').Value = '  -3.48%  '
$ws.Range('D44').Value = '2.074.03'
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '84.85'
$ws.Range('E46').Value = '  -13.71%  '
$ws.Range('E47').Value = '  +2.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.72'
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').Value = '2.791.18'
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.86'
$ws.Range('E50').Value = '  -3.96%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.64'
$ws.Range('E51').Value = '  -10.15%  '
